# Updates the "Informe-01-010010-A-TC-TM-TP" metadata sheet:
#  - row 1 headers become human readable (capitalised Spanish labels)
#  - rows 2-4 (sdmx/iaest DSD metadata: measure/dimension, medida/dim, type/URI)
#    are re-aligned per-column now that "viviendas-en-el-edificio" (column B)
#    carries its own dimension metadata instead of reusing "null"
#  - the "mapping-viviendas-en-el-edificio.xlsx" reference moves from F5 to B5
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing cell format (style index 1) onto B5 before moving the value there
$ws.Range("F5").Copy()
$ws.Range("B5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A1").Value = 'Número de edificios'
$ws.Range("B1").Value = 'Viviendas en el edificio'
$ws.Range("C1").Value = 'Comarca nombre'
$ws.Range("D1").Value = 'Comarca código'
$ws.Range("E1").Value = 'Provincia código'
$ws.Range("F1").Value = 'Aragón'
$ws.Range("G1").Value = 'Municipio código'
$ws.Range("H1").Value = 'Provincia nombre'
$ws.Range("I1").Value = 'Municipio nombre'
$ws.Range("A2").Value = 'iaest-measure:numero-de-edificios'
$ws.Range("B2").Value = 'iaest-dimension:viviendas-en-el-edificio'
$ws.Range("C2").Value = 'sdmx-dimension:refArea'
$ws.Range("D2").Value = 'null'
$ws.Range("E2").Value = 'null'
$ws.Range("F2").Value = 'sdmx-dimension:refArea'
$ws.Range("G2").Value = 'null'
$ws.Range("H2").Value = 'sdmx-dimension:refArea'
$ws.Range("I2").Value = 'sdmx-dimension:refArea'
$ws.Range("A3").Value = 'medida'
$ws.Range("B3").Value = 'dim'
$ws.Range("C3").Value = 'dim'
$ws.Range("D3").Value = 'null'
$ws.Range("E3").Value = 'null'
$ws.Range("F3").Value = 'dim'
$ws.Range("G3").Value = 'null'
$ws.Range("H3").Value = 'dim'
$ws.Range("I3").Value = 'dim'
$ws.Range("A4").Value = 'xsd:int'
$ws.Range("B4").Value = 'skos:Concept'
$ws.Range("C4").Value = 'URI-comarca'
$ws.Range("D4").Value = 'null'
$ws.Range("E4").Value = 'null'
$ws.Range("F4").Value = 'URI-Comunidad'
$ws.Range("G4").Value = 'null'
$ws.Range("H4").Value = 'URI-Provincia'
$ws.Range("I4").Value = 'URI-Municipio'
$ws.Range("A5").Clear()
$ws.Range("B5").Value = 'mapping-viviendas-en-el-edificio.xlsx'
$ws.Range("C5").Clear()
$ws.Range("D5").Clear()
$ws.Range("E5").Clear()
$ws.Range("F5").Clear()  # old value moved to B5
$ws.Range("G5").Clear()
$ws.Range("H5").Clear()
$ws.Range("I5").Clear()
